$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 128 (run_id 127)
$ws.Cells.Item(128, 1).Value = 127
$ws.Cells.Item(128, 2).Value = 1
$ws.Cells.Item(128, 3).Value = "2024-06-17 15:14:33"
$ws.Cells.Item(128, 4).Value = 200
$ws.Cells.Item(128, 5).Value = 16

# New row 129 (run_id 128)
$ws.Cells.Item(129, 1).Value = 128
$ws.Cells.Item(129, 2).Value = 2
$ws.Cells.Item(129, 3).Value = "2024-06-17 15:14:34"
$ws.Cells.Item(129, 4).Value = 200
$ws.Cells.Item(129, 5).Value = 2
